# ---------------------------------------------------------------------------
# Applies the "Add files via upload" edit to guide41_movie.xlsx:
#   * inserts a new "はじめに" (Introduction) <h3> anchor heading row near the
#     top of both the "p1" and "p2" sheets (sheet2.xml / sheet3.xml)
#   * on "p2" also inserts the existing explanatory paragraph (shared string
#     already used once on that sheet) right under the new heading
#   * makes "p2" the active sheet/tab and updates the remembered selections
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # "p1"
$ws3 = $wb.Worksheets.Item(3)   # "p2"

# ---------------------------------------------------------------------------
# 1) Sheet "p1" (sheet2.xml): insert a new row at position 7.
#    Old row 7 (blank, style s=11) shifts to row 8; old row 8 (value 35)
#    shifts to row 9; everything below shifts by one as well.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(7).Insert(-4121)   # -4121 = xlShiftDown

# After the insert:
#   row5 still holds the old row5 content (value 82, ht 207.75)
#   row6 is untouched
#   row7 is a blank placeholder row
#   row8 holds the old row7 content (blank, style s=11)
#   row9 holds the old row8 content (value 35) -- already correct, untouched

# Move row5's content (the long "82" paragraph + its row height) down to row8,
# where it belongs after the insert.
$ws2.Range("A5:B5").Copy()
$ws2.Range("A8").PasteSpecial(-4104)   # -4104 = xlPasteAll
$ws2.Rows.Item(8).RowHeight = $ws2.Rows.Item(5).RowHeight

# Row5 becomes the (now blank) row that used to be row7: same style as B4
# (style s=11), no value, default row height.
$ws2.Range("B4").Copy()
$ws2.Range("B5").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$ws2.Range("A5:B5").ClearContents()
$ws2.Rows.Item(5).RowHeight = $ws2.Rows.Item(6).RowHeight

# Row7 becomes the new "はじめに" heading row: style copied from B9 (style
# s=12), 18.75pt tall.
$ws2.Range("B9").Copy()
$ws2.Range("B7").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$intro = [char]0x003C + "h3" + [char]0x003E + [char]0x003C + "a name=" + [char]0x0022 + "first" + [char]0x0022 + [char]0x003E + [char]0x003C + "/a" + [char]0x003E + [char]0x306F + [char]0x3058 + [char]0x3081 + [char]0x306B + [char]0x003C + "/h3" + [char]0x003E
$ws2.Range("B7").Value = $intro
$ws2.Range("B7").Characters(14, 5).Font.Name = "游ゴシック"
$ws2.Range("B7").Characters(14, 5).Font.Size = 11
$ws2.Range("B7").Characters(14, 5).Font.Color = 0
$ws2.Range("B7").Characters(19, 15).Font.Name = "ＭＳ Ｐゴシック "
$ws2.Range("B7").Characters(19, 15).Font.Size = 11
$ws2.Range("B7").Characters(19, 15).Font.Color = 0
$ws2.Rows.Item(7).RowHeight = 18.75

# ---------------------------------------------------------------------------
# 2) Sheet "p2" (sheet3.xml): insert two new rows at position 7, for the new
#    heading + the accompanying paragraph (re-using the shared string that
#    row1 already uses).
# ---------------------------------------------------------------------------
$ws3.Rows.Item(7).Resize(2).Insert(-4121)   # -4121 = xlShiftDown

# New row7: same "はじめに" heading as p1!B7 -- copy it wholesale so both
# cells end up pointing at the very same shared string.
$ws2.Range("B7").Copy()
$ws3.Range("B7").PasteSpecial(-4104)   # -4104 = xlPasteAll
$ws3.Rows.Item(7).RowHeight = 18.75

# New row8: the explanatory paragraph already used by p2!B1 -- copy it so the
# new cell reuses that same shared string instead of creating a duplicate.
$ws3.Range("B1").Copy()
$ws3.Range("B8").PasteSpecial(-4104)   # -4104 = xlPasteAll
$ws3.Rows.Item(8).RowHeight = $ws3.Rows.Item(1).RowHeight

# ---------------------------------------------------------------------------
# 3) Activate "p2" and restore the remembered selections on both sheets.
# ---------------------------------------------------------------------------
$ws2.Range("B7").Select()
$ws3.Activate()
$ws3.Range("B12").Select()
